# edit.ps1 - apply the Readme.docx OOXML diff via Word COM-interop.
#
# Strategy: for each affected paragraph, replace its whole Range (start..end,
# including the trailing paragraph mark) with a hand-built <w:p> fragment via
# Range.InsertXML. When the replacement XML omits <w:pPr>, the host preserves
# the paragraph's existing <w:pPr> (and the <w:p> attributes such as
# w14:paraId / rsids) untouched - that is how "in place" run restructuring
# is done below. Brand-new paragraphs include an explicit <w:pPr> so they
# match what a freshly-typed Word paragraph looks like (no paraId/rsids).
#
# Edits are applied from the bottom of the document upward so that
# paragraph indices used for edits still to come are never invalidated by
# paragraphs inserted earlier in this script.

$d = $word.ActiveDocument
$wns = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

function Replace-ParagraphRuns($paraIndex, $innerXml) {
    $p = $d.Paragraphs($paraIndex)
    $full = $p.Range
    $target = $d.Range($full.Start, $full.End)
    $xml = "<w:p $wns>" + $innerXml + "</w:p>"
    $target.InsertXML($xml)
}

function Insert-NewParagraphAfter($paraIndex, $innerXml) {
    $p = $d.Paragraphs($paraIndex)
    $full = $p.Range
    $anchor = $d.Range($full.Start, $full.End)
    $anchor.InsertParagraphAfter()
    $nextIndex = $paraIndex + 1
    $newp = $d.Paragraphs($nextIndex)
    $newrange = $newp.Range
    $target = $d.Range($newrange.Start, $newrange.End)
    $xml = "<w:p $wns>" + $innerXml + "</w:p>"
    $target.InsertXML($xml)
}

# --- Paragraph 31: "Once the app is running:" -> drop <w:lastRenderedPageBreak/> ---
$inner31 = '<w:r><w:t>Once the app is running:</w:t></w:r>'
Replace-ParagraphRuns 31 $inner31

# --- Insert a new empty (bold) paragraph between the 3rd blank paragraph (29) and "Usage" (30) ---
$innerBlank = '<w:pPr><w:rPr><w:b/><w:bCs/></w:rPr></w:pPr>'
Insert-NewParagraphAfter 29 $innerBlank

# --- Paragraph 25: "streamlit run app.py" (bold) -> split "streamlit" out with proofErr ---
$inner25 = '<w:proofErr w:type="spellStart"/>' + `
    '<w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>streamlit</w:t></w:r>' + `
    '<w:proofErr w:type="spellEnd"/>' + `
    '<w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t xml:space="preserve"> run app.py</w:t></w:r>'
Replace-ParagraphRuns 25 $inner25

# --- Paragraph 24: "Run the Streamlit app using the following command:" ---
$inner24 = '<w:r><w:t xml:space="preserve">Run the </w:t></w:r>' + `
    '<w:proofErr w:type="spellStart"/>' + `
    '<w:r><w:t>Streamlit</w:t></w:r>' + `
    '<w:proofErr w:type="spellEnd"/>' + `
    '<w:r><w:t xml:space="preserve"> app using the following command:</w:t></w:r>'
Replace-ParagraphRuns 24 $inner24

# --- Paragraph 23: "5. Run the Streamlit App" (bold) ---
$inner23 = '<w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t xml:space="preserve">5. Run the </w:t></w:r>' + `
    '<w:proofErr w:type="spellStart"/>' + `
    '<w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>Streamlit</w:t></w:r>' + `
    '<w:proofErr w:type="spellEnd"/>' + `
    '<w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t xml:space="preserve"> App</w:t></w:r>'
Replace-ParagraphRuns 23 $inner23

# --- New paragraph inserted right after paragraph 22 (old "pip install ..." text moved here, split up) ---
$innerNewPip = '<w:pPr><w:rPr><w:b/><w:bCs/></w:rPr></w:pPr>' + `
    '<w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>pip install google-</w:t></w:r>' + `
    '<w:proofErr w:type="spellStart"/>' + `
    '<w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>generativeai</w:t></w:r>' + `
    '<w:proofErr w:type="spellEnd"/>' + `
    '<w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t xml:space="preserve"> </w:t></w:r>' + `
    '<w:proofErr w:type="spellStart"/>' + `
    '<w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>streamlit</w:t></w:r>' + `
    '<w:proofErr w:type="spellEnd"/>'
Insert-NewParagraphAfter 22 $innerNewPip

# --- Paragraph 22: "pip install google-generativeai streamlit" -> "pip install -r requirements.txt" (bold) ---
$inner22 = '<w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>pip install -r requirements.txt</w:t></w:r>'
Replace-ParagraphRuns 22 $inner22

# --- Paragraph 19: "cd DoJ-Chatbot-using-Generative-AI" (bold) ---
$inner19 = '<w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t xml:space="preserve">cd </w:t></w:r>' + `
    '<w:proofErr w:type="spellStart"/>' + `
    '<w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>DoJ</w:t></w:r>' + `
    '<w:proofErr w:type="spellEnd"/>' + `
    '<w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>-Chatbot-using-Generative-AI</w:t></w:r>'
Replace-ParagraphRuns 19 $inner19

# --- Paragraph 14: "Required libraries: google-generativeai, streamlit" ---
$inner14 = '<w:r><w:t>Required libraries: google-</w:t></w:r>' + `
    '<w:proofErr w:type="spellStart"/>' + `
    '<w:r><w:t>generativeai</w:t></w:r>' + `
    '<w:proofErr w:type="spellEnd"/>' + `
    '<w:r><w:t xml:space="preserve">, </w:t></w:r>' + `
    '<w:proofErr w:type="spellStart"/>' + `
    '<w:r><w:t>streamlit</w:t></w:r>' + `
    '<w:proofErr w:type="spellEnd"/>'
Replace-ParagraphRuns 14 $inner14

# --- Paragraph 13: "Streamlit installed" ---
$inner13 = '<w:proofErr w:type="spellStart"/>' + `
    '<w:r><w:t>Streamlit</w:t></w:r>' + `
    '<w:proofErr w:type="spellEnd"/>' + `
    '<w:r><w:t xml:space="preserve"> installed</w:t></w:r>'
Replace-ParagraphRuns 13 $inner13

# --- Paragraph 6: "Styled UI using Streamlit for easy interaction." ---
$inner6 = '<w:r><w:t xml:space="preserve">Styled UI using </w:t></w:r>' + `
    '<w:proofErr w:type="spellStart"/>' + `
    '<w:r><w:t>Streamlit</w:t></w:r>' + `
    '<w:proofErr w:type="spellEnd"/>' + `
    '<w:r><w:t xml:space="preserve"> for easy interaction.</w:t></w:r>'
Replace-ParagraphRuns 6 $inner6

# --- Paragraph 2: "PravdaGPT is a Generative AI-powered chatbot ..." ---
$inner2 = '<w:proofErr w:type="spellStart"/>' + `
    '<w:r><w:t>PravdaGPT</w:t></w:r>' + `
    '<w:proofErr w:type="spellEnd"/>' + `
    '<w:r><w:t xml:space="preserve"> is a Generative AI-powered chatbot designed to assist with legal and Department of Justice (DOJ) related queries. It provides answers to questions related to penal codes, DOJ policies, legal terms, and more.</w:t></w:r>'
Replace-ParagraphRuns 2 $inner2

# --- Paragraph 1: "PravdaGPT - Legal and DOJ Specialist Assistant" (bold) ---
$inner1 = '<w:proofErr w:type="spellStart"/>' + `
    '<w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>PravdaGPT</w:t></w:r>' + `
    '<w:proofErr w:type="spellEnd"/>' + `
    '<w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t xml:space="preserve"> - Legal and DOJ Specialist Assistant</w:t></w:r>'
Replace-ParagraphRuns 1 $inner1

Write-Output "edits applied"
